$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 13167989
$ws.Range("I62").Value = 25013458
$ws.Range("J62").Value = 6355.5557
$ws.Range("K62").Value = 25013458
$ws.Range("L62").Value = 6355.5557
$ws.Range("M62").Value = -25012834
$ws.Range("N62").Value = -7603.5557
$ws.Range("H65").Value = 13167989
$ws.Range("I65").Value = 25013458
$ws.Range("J65").Value = 6355.5557
$ws.Range("K65").Value = 125067290
$ws.Range("L65").Value = 31777.7785
$ws.Range("M65").Value = -125064170
$ws.Range("N65").Value = -38017.7785
$ws.Range("H112").Value = 20736.191
$ws.Range("I112").Value = 500
$ws.Range("J112").Value = 21229.756
$ws.Range("K112").Value = 1500
$ws.Range("L112").Value = 63689.268
$ws.Range("M112").Value = -392
$ws.Range("N112").Value = -65905.26800000001
$ws.Range("H131").Value = 4063.889
$ws.Range("I131").Value = 995
$ws.Range("J131").Value = 5598.3335
$ws.Range("K131").Value = 2985
$ws.Range("L131").Value = 16795.0005
$ws.Range("M131").Value = 2055
$ws.Range("N131").Value = -26875.0005
$ws.Range("H132").Value = 21981540
$ws.Range("I132").Value = 18075310
$ws.Range("J132").Value = 62508676
$ws.Range("K132").Value = 54225930
$ws.Range("L132").Value = 187526028
$ws.Range("M132").Value = -54223400
$ws.Range("N132").Value = -187531088
$ws.Range("H137").Value = 1439.238
$ws.Range("I137").Value = 1202.8667
$ws.Range("J137").Value = 1570.5555
$ws.Range("K137").Value = 3608.6001
$ws.Range("L137").Value = 4711.666499999999
$ws.Range("M137").Value = -1058.6001
$ws.Range("N137").Value = -9811.666499999999
$ws.Range("H138").Value = 3163.3958
$ws.Range("I138").Value = 1597.7567
$ws.Range("J138").Value = 4145.2373
$ws.Range("K138").Value = 4793.2701
$ws.Range("L138").Value = 12435.7119
$ws.Range("M138").Value = 346.7299000000003
$ws.Range("N138").Value = -22715.7119

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12093.333
$ws.Range("I32").Value = 11065.031
$ws.Range("J32").Value = 44999
$ws.Range("K32").Value = 11065.031
$ws.Range("L32").Value = 44999
$ws.Range("M32").Value = -10778.031
$ws.Range("N32").Value = -45573
$ws.Range("H61").Value = 1708.6207
$ws.Range("I61").Value = 1608.7742
$ws.Range("J61").Value = 1823.2593
$ws.Range("K61").Value = 1608.7742
$ws.Range("L61").Value = 1823.2593
$ws.Range("M61").Value = -1396.7742
$ws.Range("N61").Value = -2247.2593
$ws.Range("H136").Value = 1708.6207
$ws.Range("I136").Value = 1608.7742
$ws.Range("J136").Value = 1823.2593
$ws.Range("K136").Value = 4826.3226
$ws.Range("L136").Value = 5469.7779
$ws.Range("M136").Value = -2276.3226
$ws.Range("N136").Value = -10569.7779

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 575.35486
$ws.Range("I94").Value = 465.27274
$ws.Range("J94").Value = 844.44446
$ws.Range("K94").Value = 465.27274
$ws.Range("L94").Value = 844.44446
$ws.Range("M94").Value = -14.27274
$ws.Range("N94").Value = -1746.44446
$ws.Range("H132").Value = 45544.445
$ws.Range("J132").Value = 45544.445
$ws.Range("L132").Value = 45544.445
$ws.Range("N132").Value = -55664.445

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 4250
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").Value = ""
$ws.Range("H69").Value = 5000
$ws.Range("I69").Value = 5000
$ws.Range("K69").Value = 5000
$ws.Range("M69").Value = -4251
$ws.Range("H72").Value = 5000
$ws.Range("I72").Value = 5000
$ws.Range("K72").Value = 15000
$ws.Range("M72").Value = -11256
$ws.Range("H87").Value = 9800
$ws.Range("J87").Value = 9800
$ws.Range("L87").Value = 9800
$ws.Range("N87").Value = -12172
$ws.Range("H90").Value = 9800
$ws.Range("J90").Value = 9800
$ws.Range("L90").Value = 29400
$ws.Range("N90").Value = -41256
$ws.Range("H99").Value = 52639336
$ws.Range("I99").Value = 90920024
$ws.Range("J99").Value = 3387.5
$ws.Range("K99").Value = 90920024
$ws.Range("L99").Value = 3387.5
$ws.Range("M99").Value = -90918526
$ws.Range("N99").Value = -6383.5
$ws.Range("H126").Value = 52639336
$ws.Range("I126").Value = 90920024
$ws.Range("J126").Value = 3387.5
$ws.Range("K126").Value = 272760072
$ws.Range("L126").Value = 10162.5
$ws.Range("M126").Value = -272757602
$ws.Range("N126").Value = -15102.5
$ws.Range("H141").Value = 56345.453
$ws.Range("I141").Value = 16500
$ws.Range("J141").Value = 65200
$ws.Range("K141").Value = 16500
$ws.Range("L141").Value = 65200
$ws.Range("M141").Value = -11320
$ws.Range("N141").Value = -75560

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 398.33334
$ws.Range("I13").Value = 398.33334
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 1195.00002
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -1027.00002
$ws.Range("N13").Value = ""
$ws.Range("H131").Value = 732.5306399999999
$ws.Range("I131").Value = 438.75
$ws.Range("J131").Value = 758.6445
$ws.Range("K131").Value = 1316.25
$ws.Range("L131").Value = 2275.9335
$ws.Range("M131").Value = 3723.75
$ws.Range("N131").Value = -12355.9335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = ""
$ws.Range("H132").Value = 5038.2256
$ws.Range("I132").Value = 1258.5
$ws.Range("J132").Value = 9069.933999999999
$ws.Range("K132").Value = 3775.5
$ws.Range("L132").Value = 27209.802
$ws.Range("M132").Value = -1245.5
$ws.Range("N132").Value = -32269.802

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2194.4666
$ws.Range("I7").Value = 2129.7144
$ws.Range("K7").Value = 2129.7144
$ws.Range("M7").Value = -2017.7144
$ws.Range("H40").Value = 17859094
$ws.Range("I40").Value = 1781.2
$ws.Range("J40").Value = 62502376
$ws.Range("K40").Value = 1781.2
$ws.Range("L40").Value = 62502376
$ws.Range("M40").Value = -1645.2
$ws.Range("N40").Value = -62502648
$ws.Range("H93").Value = 1247
$ws.Range("I93").Value = 1213.3572
$ws.Range("J93").Value = 1325.5
$ws.Range("K93").Value = 1213.3572
$ws.Range("L93").Value = 1325.5
$ws.Range("M93").Value = 34.64280000000008
$ws.Range("N93").Value = -3821.5
$ws.Range("H109").Value = 34541.25
$ws.Range("J109").Value = 34541.25
$ws.Range("L109").Value = 34541.25
$ws.Range("N109").Value = -37315.25
$ws.Range("H126").Value = 2194.4666
$ws.Range("I126").Value = 2129.7144
$ws.Range("K126").Value = 6389.1432
$ws.Range("M126").Value = -3919.1432
$ws.Range("H136").Value = 4159.1025
$ws.Range("I136").Value = 2871.4285
$ws.Range("J136").Value = 7436.8184
$ws.Range("K136").Value = 8614.2855
$ws.Range("L136").Value = 22310.4552
$ws.Range("M136").Value = -6064.2855
$ws.Range("N136").Value = -27410.4552

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H117").Value = 25000
$ws.Range("J117").Value = 25000
$ws.Range("L117").Value = 25000
$ws.Range("N117").Value = -34178
$ws.Range("H122").Value = 1971.6897
$ws.Range("I122").Value = 1256.1428
$ws.Range("J122").Value = 3850
$ws.Range("K122").Value = 3768.4284
$ws.Range("L122").Value = 11550
$ws.Range("M122").Value = -1318.4284
$ws.Range("N122").Value = -16450
$ws.Range("H126").Value = 1337.35
$ws.Range("I126").Value = 967.6429000000001
$ws.Range("K126").Value = 2902.9287
$ws.Range("M126").Value = -432.9287000000004
